$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 369, shifting existing rows 369:379 down to 370:380
$ws.Range("A369").EntireRow.Insert()

# Populate the newly inserted row 369 with data
$ws.Range("A369").Value = 10
$ws.Range("B369").Value = "Vega Modelo de Temuco"
$ws.Range("C369").Value = "La Araucanía"
$ws.Range("D369").Value = 45075
$ws.Range("E369").Value = 9
$ws.Range("F369").Value = "Fruta"
$ws.Range("G369").Value = 100103
$ws.Range("H369").Value = "Frutos de hueso (carozo)"
$ws.Range("I369").Value = 100103002
$ws.Range("J369").Value = "Ciruela"
$ws.Range("K369").Value = "Blue Giant"
$ws.Range("L369").Value = "Primera"
$ws.Range("M369").Value = 150
$ws.Range("N369").Value = 15000
$ws.Range("O369").Value = 15000
$ws.Range("P369").Value = 15000
$ws.Range("Q369").Value = "`$/bandeja 18 kilos granel"
$ws.Range("R369").Value = "Región de O'Higgins"
$ws.Range("S369").Value = 833
$ws.Range("T369").Value = 18
